$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting the existing rows 84:98 down to 85:99
$ws.Rows("84:84").Insert()

$row = 84
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item($row, 3).Value = 'La Araucanía'
$ws.Cells.Item($row, 4).Value = 44855
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 'Fruta'
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = 'Berries'
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = 'Arándano (blue)'
$ws.Cells.Item($row, 11).Value = 'Sin especificar'
$ws.Cells.Item($row, 12).Value = 'Primera'
$ws.Cells.Item($row, 13).Value = 25
$ws.Cells.Item($row, 14).Value = 7000
$ws.Cells.Item($row, 15).Value = 7000
$ws.Cells.Item($row, 16).Value = 7000
$ws.Cells.Item($row, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item($row, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item($row, 19).Value = 3500
$ws.Cells.Item($row, 20).Value = 2
